$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was re-sorted and the "Biltin" item (BSL 2 / ISL 1) was dropped,
# so the sheet now has 26 data rows (27 total including header) instead of 27 data rows.
# Remove the last row first so row indices/dimension line up, then rewrite rows 2-27
# with the new static target/sales data layout.
$ws.Rows("28:28").Delete()

$data = @(
    @('4', 'Desodin', '1', 'Desodin 60ml Syrup', '60 ml', '0', '0'),
    @('5', 'Dinafex', '2', 'Dinafex 180mg Tablet', '30''s', '0', '0'),
    @('5', 'Dinafex', '3', 'Dinafex 120mg Tablet', '30''s', '0', '0'),
    @('5', 'Dinafex', '4', 'Dinafex 60mg Tablet', '30''s', '0', '0'),
    @('6', 'Dorenta', '5', 'Dorenta 50mg Tablet', '50''s', '0', '0'),
    @('7', 'Etorix', '6', 'Etorix 60mg Tablet - 40''s', '40''s', '0', '0'),
    @('7', 'Etorix', '7', 'Etorix 120mg Tablet', '20''s', '0', '0'),
    @('7', 'Etorix', '8', 'Etorix 90mg Tablet', '30''s', '0', '0'),
    @('8', 'Fenobac', '9', 'Fenobac 100ml Syrup', '100ml', '0', '0'),
    @('9', 'Flucloxin', '10', 'Flucloxin 500mg Capsule', '30 ''s', '0', '0'),
    @('9', 'Flucloxin', '11', 'Flucloxin 500mg Capsule - 36''s', '36 ''s', '0', '0'),
    @('10', 'Geminox', '12', 'Geminox 320mg Tablet - 8''s', '8 ''s', '0', '0'),
    @('11', 'Ketonic', '13', 'Ketonic 30mg Injection', '5 ''s', '0', '0'),
    @('11', 'Ketonic', '14', 'Ketonic 30mg IM/IV Injection - 4''s', '4''s', '0', '0'),
    @('11', 'Ketonic', '15', 'Ketonic 10mg Tablet', '20''s', '0', '0'),
    @('12', 'Kynol', '16', 'Kynol TR 200mg Capsule', '30 ''s', '0', '0'),
    @('12', 'Kynol', '17', 'Kynol D 25mg Tablet', '60 ''s', '0', '0'),
    @('12', 'Kynol', '18', 'Kynol TR 100mg Capsule', '50 ''s', '0', '0'),
    @('17', 'Naprox', '19', 'Naprox Plus 500mg Tablet - 30''s', '30 ''s', '0', '0'),
    @('19', 'Oradin', '20', 'Oradin Plus Tablet - 40''s', '40 ''s', '0', '0'),
    @('20', 'Osticare', '21', 'Osticare Tablet 24''s', '24''s', '0', '0'),
    @('23', 'Rupaday', '22', 'Rupaday Oral Solution 60ml', '1''s', '0', '0'),
    @('35', 'Zithrox', '23', 'Zithrox 30ml Dry Suspension', '30ml', '0', '0'),
    @('35', 'Zithrox', '24', 'Zithrox 250mg Tablet - 6''s', '6''s', '0', '0'),
    @('35', 'Zithrox', '25', 'Zithrox 500mg Tablet', '6 ''s', '0', '0'),
    @('35', 'Zithrox', '26', 'Zithrox 15ml Suspension', '15 ml', '0', '0')
)

$cols = @("A", "B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $vals[$j]
    }
}
